{"js": "// Remove the \"Classification: Controlled\" footer content-marking text boxes\n// (floating shapes) from every section's footers. These shapes were added by\n// the classification/sensitivity-labeling feature and live in the footer's\n// shape collection rather than in the footer body's plain text, so we walk\n// Section -> Footer(primary/firstPage/evenPages) -> shapes and delete any\n// shape we find there.\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst footerTypes = [\"primary\", \"firstPage\", \"evenPages\"];\n\nfor (let i = 0; i < sections.items.length; i++) {\n  const section = sections.items[i];\n\n  for (const footerType of footerTypes) {\n    const footer = section.getFooter(footerType);\n    const shapes = footer.shapes;\n    shapes.load(\"items\");\n    await context.sync();\n\n    // Delete from the end so indices stay valid while we remove items.\n    for (let j = shapes.items.length - 1; j >= 0; j--) {\n      shapes.items[j].delete();\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the \"Classification: Controlled\" footer content-marking text boxes\n# (floating shapes) from every section's footers. These shapes were added by\n# the classification/sensitivity-labeling feature and live in the footer's\n# Shapes collection rather than in the footer Range's plain text, so we walk\n# Section -> Footers(wdHeaderFooterPrimary/FirstPage/EvenPages) -> Shapes and\n# delete anything we find there.\n\n$d = $word.ActiveDocument\n\nfor ($s = 1; $s -le $d.Sections.Count; $s++) {\n    $section = $d.Sections.Item($s)\n    $footers = $section.Footers\n\n    for ($f = 1; $f -le $footers.Count; $f++) {\n        $footer = $footers.Item($f)\n        while ($footer.Shapes.Count -gt 0) {\n            $footer.Shapes.Item(1).Delete()\n        }\n    }\n}\n"}
